$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds plain-text numbers (e.g. "26.681.72", "211.44").
# Force text format before writing so Excel does not auto-convert simple
# decimals (e.g. "211.48") into Number cells, then restore General/Normal
# so the cell style matches the original (no explicit style index).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.682.62"
$ws.Range("D3").Value = "1.598.01"
$ws.Range("D5").Value = "211.48"
$ws.Range("D10").Value = "19.48"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D12").Value = "1.822.73"
$ws.Range("D13").Value = "1.595.22"
$ws.Range("D16").Value = "65.12"
$ws.Range("D17").Value = "26.657.75"
$ws.Range("D18").Value = "0.0₃0749"
$ws.Range("D19").Value = "209.85"
$ws.Range("D21").Value = "7.04"
$ws.Range("D25").Value = "143.38"
$ws.Range("D26").Value = "1.01"
$ws.Range("D27").Value = "7.11"
$ws.Range("D29").Value = "15.33"
$ws.Range("D30").Value = "0.0517"
$ws.Range("D34").Value = "1.287.58"
$ws.Range("D40").Value = "0.826"
$ws.Range("D43").Value = "0.782"
$ws.Range("D44").Value = "63.21"
$ws.Range("D45").Value = "1.734.84"
$ws.Range("D46").Value = "91.02"
$ws.Range("D47").Value = "1.57"
$ws.Range("D48").Value = "0.101"
$ws.Range("D51").Value = "7.33"

$ws.Range("D2:D51").Style = "Normal"

# Volume(1h) column (E) holds percentage text with 2 leading/trailing spaces.
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  +17.75%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -2.37%  "
